$wb = $excel.ActiveWorkbook

# The workbook's recorded output path changed (new machine / new folder
# layout: "Katalon_mastercopy520\DataCommons_Automation" instead of
# "Users\radhakrishnang2\Desktop\Commons_Automation"). Update every cell
# that echoes the output file path to the new value.
$newPath = "C:\Katalon_mastercopy520\DataCommons_Automation\OutputFiles\TC01_Canine_Filter_Breed-Akita_Neo4jData.xlsx"

# "CypherOutput_Message" sheet: row 10 holds the "Output:" file path.
$wsCypherMsg = $wb.Worksheets.Item("CypherOutput_Message")
$wsCypherMsg.Range("A10").Value = $newPath

# "StatOutput_Message" sheet: the same block repeats twice (rows 1-10 and
# 11-20), each ending with the output path on row 10 / row 20.
$wsStatMsg = $wb.Worksheets.Item("StatOutput_Message")
$wsStatMsg.Range("A10").Value = $newPath
$wsStatMsg.Range("A20").Value = $newPath
